# Commit: "Fixed CCDI queries and ICDC index issue"
#
# The Participant-tab SQL query in cell B2 of Sheet1 is corrected:
#   - the "Ethnicity" column is dropped from the SELECT list
#   - "Alternate ID" is renamed to "Synonym Participant ID"
# All other tabs/cells are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newParticipantQuery = @"
SELECT
    p.participant_id AS "Participant ID",
    st.study_id AS "Study ID",
    COALESCE(p.sex_at_birth, '') AS "Sex",
    COALESCE(p.race, '') AS "Race",
    COALESCE(CAST(p.alternate_participant_id AS INT), '') AS "Synonym Participant ID"
FROM 
    df_participant p
JOIN 
    df_study st ON p."study.id" = st.id
LEFT JOIN 
    df_sample smp ON smp."participant.id" = p.participant_id
LEFT JOIN 
    df_diagnosis diag ON diag."participant.id" = p.participant_id
LEFT JOIN 
    df_publication pub ON pub."study.id" = st.study_id
LEFT JOIN 
    df_sequencing_file sqf ON sqf."sample.id" = smp.sample_id
WHERE 
    st.study_id = 'phs002371' AND p.sex_at_birth = 'Male'
ORDER BY 
    p.participant_id ASC
LIMIT 100;
"@

# Re-stamp the cell's font so Excel records it against a fresh style entry
# (mirrors what happened in the source edit - the cell keeps the same visual
# look: 12pt Calibri, theme color 1, wrapped text) and then write the new text.
$ws.Range("B2").Font.ThemeColor = 1
$ws.Range("B2").WrapText = $true
$ws.Range("B2").Value = $newParticipantQuery

# Leave the selection where the source workbook left it when saved.
$ws.Range("C2").Select()
